$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 2-39) currently holds date serial numbers formatted with a
# custom "YYYY-MM-DD HH:MM:SS" number format. The naive forecaster bug fix
# replaces those with plain text quarter labels ("1987Q4", "1988Q4", ...),
# and the cells should take on the same look as the header (no special
# date number format).

$labels = @(
    "1987Q4", "1988Q4", "1989Q4", "1990Q4", "1991Q4", "1992Q4", "1993Q4",
    "1994Q4", "1995Q4", "1996Q4", "1997Q4", "1998Q4", "1999Q4", "2000Q4",
    "2001Q4", "2002Q4", "2003Q4", "2004Q4", "2005Q4", "2006Q4", "2007Q4",
    "2008Q4", "2009Q4", "2010Q4", "2011Q4", "2012Q4", "2013Q4", "2014Q4",
    "2015Q4", "2016Q4", "2017Q4", "2018Q4", "2019Q4", "2020Q4", "2021Q4",
    "2022Q4", "2023Q4", "2024Q4"
)

$ws.Range("A1").Copy()

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.PasteSpecial(-4122)
    $cell.Value = $labels[$i]
}

$excel.CutCopyMode = $false
